# High voltage CT (current transformer) BOM additions, hyperlinks, and view updates.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new BOM rows (30 and 31) for the CR Magnetics current transformers ---
# Copy formatting (style) from the last existing data row (29) down onto the
# two new rows so they pick up the same cell style (s="1") without creating
# redundant style entries.
$ws.Range("A29:G29").Copy() | Out-Null
$ws.Range("A30:G30").PasteSpecial(-4122) | Out-Null
$ws.Range("A31:G31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 30: 100A current transformer
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "100 A current Transformer"
$ws.Range("C30").Value = ""
$ws.Range("D30").Value = "Curr Sense Xfmr 100A In-line 1000:1"
$ws.Range("E30").Value = "CR Magnetics"
$ws.Range("F30").Value = "CR8450-1000"
$ws.Range("G30").Value = "https://octopart.com/cr8450-1000-cr+magnetics-19500880?r=sp"

# Row 31: 4A current transformer (optional)
$ws.Range("A31").Value = "optional"
$ws.Range("B31").Value = "4 A current Transformer"
$ws.Range("C31").Value = ""
$ws.Range("D31").Value = "Curr Sense Xfmr 4A In-line 1000:1"
$ws.Range("E31").Value = "CR Magnetics"
$ws.Range("F31").Value = "CR8401-1000-G"
$ws.Range("G31").Value = "https://octopart.com/cr8401-1000-g-cr+magnetics-1366028?r=sp"

# --- Turn the OCTOPART_URL cells into real hyperlinks (adds the "Hyperlink" ---
# --- cell style/font, matching the order rId1..rId4 seen in the target file) ---
$ws.Hyperlinks.Add($ws.Range("G30"), "https://octopart.com/cr8450-1000-cr+magnetics-19500880?r=sp") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G20"), "https://octopart.com/search?q=RC0603JR-0710KL&start=0") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G21"), "https://octopart.com/crcw06031r00fkea-vishay-39713722?r=sp&s=biXSlGylTBiwPOCcABx_kA") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G31"), "https://octopart.com/cr8401-1000-g-cr+magnetics-1366028?r=sp") | Out-Null

# --- Page setup: switch to portrait orientation (adds pageSetup element) ---
$ws.PageSetup.Orientation = 1

# --- Update view state: scroll down and leave selection on B36 ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B36").Select() | Out-Null

Write-Host "Applied CT BOM additions, hyperlinks, and page setup changes."
